# Update countries & provincias Spain
# - Armenia overtakes Austria in total-cases ranking (rows 60/61 swap
#   underlying data while keeping Austria's figures unchanged).
# - Israel (row 28), Argelia (row 71) and Republica de Macedonia (row 94)
#   get refreshed daily figures.
# - The "last updated" timestamp moves from 08:22 to 09:39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name swap: Armenia now ranks above Austria ---------------
$ws.Range("A60").Value = "Armenia"
$ws.Range("A61").Value = "Austria"

# --- Row 60 (now Armenia): updated case counts ------------------------
$ws.Range("B60").Value = 66694
$ws.Range("C60").Value = 1234
$ws.Range("D60").Value = 48734
$ws.Range("E60").Value = 16859
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 10
$ws.Range("H60").Value = 1101

# --- Row 61 (now Austria): keeps Austria's previous figures -----------
$ws.Range("B61").Value = 65927
$ws.Range("C61").Value = 0
$ws.Range("D61").Value = 50359
$ws.Range("E61").Value = 14664
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 904

# --- Row 28 (Israel): refreshed figures --------------------------------
$ws.Range("B28").Value = 305348
$ws.Range("C28").Value = 472
$ws.Range("D28").Value = 279729
$ws.Range("E28").Value = 23351
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 2268

# --- Row 71 (Argelia): refreshed figures -------------------------------
$ws.Range("B71").Value = 48757
$ws.Range("C71").Value = 989
$ws.Range("D71").Value = 14637
$ws.Range("E71").Value = 32909
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 38
$ws.Range("H71").Value = 1211

# --- Row 94 (Republica de Macedonia): refreshed figures ----------------
$ws.Range("B94").Value = 19857
$ws.Range("C94").Value = 1194
$ws.Range("D94").Value = 8666
$ws.Range("E94").Value = 11033
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 15
$ws.Range("H94").Value = 158

# --- Timestamp footer row ----------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 20 de Octubre de 2020 a las 09:39"
